# Updated run for publication
# Refresh frequency-table values in rows 2-5 (bases A, C, G, T), columns B:X (positions 1-23)
# with the newly computed per-position frequencies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.00681302043906132
$ws.Cells.Item(2, 3).Value = 0.00151400454201363
$ws.Cells.Item(2, 4).Value = 0.00151400454201363
$ws.Cells.Item(2, 5).Value = 0.000757002271006813
$ws.Cells.Item(2, 6).Value = 0.000757002271006813
$ws.Cells.Item(2, 7).Value = 0.00151400454201363
$ws.Cells.Item(2, 8).Value = 0.0416351249053747
$ws.Cells.Item(2, 9).Value = 0.0355791067373202
$ws.Cells.Item(2, 10).Value = 0.0105980317940954
$ws.Cells.Item(2, 11).Value = 0.0227100681302044
$ws.Cells.Item(2, 12).Value = 0.012112036336109
$ws.Cells.Item(2, 13).Value = 0.00454201362604088
$ws.Cells.Item(2, 14).Value = 0.00529901589704769
$ws.Cells.Item(2, 15).Value = 0.994700984102952
$ws.Cells.Item(2, 16).Value = 0.000757002271006813
$ws.Cells.Item(2, 17).Value = 0
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(2, 20).Value = 0
$ws.Cells.Item(2, 21).Value = 0.00151400454201363
$ws.Cells.Item(2, 22).Value = 0.934140802422407
$ws.Cells.Item(2, 23).Value = 0.0151400454201363
$ws.Cells.Item(2, 24).Value = 0.0060560181680545
$ws.Cells.Item(3, 2).Value = 0.987887963663891
$ws.Cells.Item(3, 3).Value = 0.00757002271006813
$ws.Cells.Item(3, 4).Value = 0.0060560181680545
$ws.Cells.Item(3, 5).Value = 0.00681302043906132
$ws.Cells.Item(3, 6).Value = 0.989401968205905
$ws.Cells.Item(3, 7).Value = 0.98107494322483
$ws.Cells.Item(3, 8).Value = 0.0060560181680545
$ws.Cells.Item(3, 9).Value = 0.945495836487509
$ws.Cells.Item(3, 10).Value = 0.0280090840272521
$ws.Cells.Item(3, 11).Value = 0.000757002271006813
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 0.000757002271006813
$ws.Cells.Item(3, 14).Value = 0
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = 0.998485995457986
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 0.0060560181680545
$ws.Cells.Item(3, 19).Value = 0.00151400454201363
$ws.Cells.Item(3, 20).Value = 0.999242997728993
$ws.Cells.Item(3, 21).Value = 0.993186979560939
$ws.Cells.Item(3, 22).Value = 0
$ws.Cells.Item(3, 23).Value = 0.00908402725208176
$ws.Cells.Item(3, 24).Value = 0
$ws.Cells.Item(4, 2).Value = 0.000757002271006813
$ws.Cells.Item(4, 3).Value = 0.00529901589704769
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0.00227100681302044
$ws.Cells.Item(4, 6).Value = 0.000757002271006813
$ws.Cells.Item(4, 7).Value = 0.00529901589704769
$ws.Cells.Item(4, 8).Value = 0.94776684330053
$ws.Cells.Item(4, 9).Value = 0.00302800908402725
$ws.Cells.Item(4, 10).Value = 0.000757002271006813
$ws.Cells.Item(4, 11).Value = 0.9666919000757
$ws.Cells.Item(4, 12).Value = 0.987130961392884
$ws.Cells.Item(4, 13).Value = 0.993943981831946
$ws.Cells.Item(4, 14).Value = 0.993943981831946
$ws.Cells.Item(4, 15).Value = 0.000757002271006813
$ws.Cells.Item(4, 16).Value = 0.000757002271006813
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = 0
$ws.Cells.Item(4, 19).Value = 0
$ws.Cells.Item(4, 20).Value = 0
$ws.Cells.Item(4, 21).Value = 0.000757002271006813
$ws.Cells.Item(4, 22).Value = 0.0643451930355791
$ws.Cells.Item(4, 23).Value = 0.9666919000757
$ws.Cells.Item(4, 24).Value = 0.98107494322483
$ws.Cells.Item(5, 2).Value = 0.00378501135503407
$ws.Cells.Item(5, 3).Value = 0.985616956850871
$ws.Cells.Item(5, 4).Value = 0.992429977289932
$ws.Cells.Item(5, 5).Value = 0.990158970476911
$ws.Cells.Item(5, 6).Value = 0.00908402725208176
$ws.Cells.Item(5, 7).Value = 0.012112036336109
$ws.Cells.Item(5, 8).Value = 0.00454201362604088
$ws.Cells.Item(5, 9).Value = 0.0158970476911431
$ws.Cells.Item(5, 10).Value = 0.960635881907646
$ws.Cells.Item(5, 11).Value = 0.00908402725208176
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = 0
$ws.Cells.Item(5, 14).Value = 0.000757002271006813
$ws.Cells.Item(5, 15).Value = 0.00454201362604088
$ws.Cells.Item(5, 16).Value = 0
$ws.Cells.Item(5, 17).Value = 1
$ws.Cells.Item(5, 18).Value = 0.993943981831946
$ws.Cells.Item(5, 19).Value = 0.998485995457986
$ws.Cells.Item(5, 20).Value = 0.000757002271006813
$ws.Cells.Item(5, 21).Value = 0.00454201362604088
$ws.Cells.Item(5, 22).Value = 0.00151400454201363
$ws.Cells.Item(5, 23).Value = 0.00908402725208176
$ws.Cells.Item(5, 24).Value = 0.0128690386071158
